$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 47
$ws.Range("I5").Value = 46
$ws.Range("K5").Value = 46
$ws.Range("M5").Value = 69
$ws.Range("H15").Value = 1920.3889
$ws.Range("I15").Value = 1920.3889
$ws.Range("K15").Value = 5761.1667
$ws.Range("M15").Value = -5592.1667
$ws.Range("H76").Value = 5964.5386
$ws.Range("J76").Value = 5964.5386
$ws.Range("L76").Value = 5964.5386
$ws.Range("N76").Value = -6594.5386
$ws.Range("H79").Value = 5964.5386
$ws.Range("J79").Value = 5964.5386
$ws.Range("L79").Value = 5964.5386
$ws.Range("N79").Value = -8148.5386
$ws.Range("H97").Value = 1911.6
$ws.Range("I97").Value = 5000
$ws.Range("J97").Value = 1568.4445
$ws.Range("K97").Value = 15000
$ws.Range("L97").Value = 4705.333500000001
$ws.Range("M97").Value = -14504
$ws.Range("N97").Value = -5697.333500000001
$ws.Range("H107").Value = 10764.667
$ws.Range("I107").Value = 12473.2
$ws.Range("K107").Value = 12473.2
$ws.Range("M107").Value = -10553.2
$ws.Range("H131").Value = 3601.5557
$ws.Range("I131").Value = 2373.4285
$ws.Range("K131").Value = 7120.2855
$ws.Range("M131").Value = -2080.2855
$ws.Range("H132").Value = 4250.5884
$ws.Range("I132").Value = 4054
$ws.Range("K132").Value = 12162
$ws.Range("M132").Value = -9632
$ws.Range("H137").Value = 8624.896000000001
$ws.Range("I137").Value = 11372.1
$ws.Range("J137").Value = 2520
$ws.Range("K137").Value = 34116.3
$ws.Range("L137").Value = 7560
$ws.Range("M137").Value = -31566.3
$ws.Range("N137").Value = -12660
$ws.Range("H141").Value = 6098.353
$ws.Range("I141").Value = 5854.5
$ws.Range("K141").Value = 17563.5
$ws.Range("M141").Value = -12383.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6728.618
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10574
$ws.Range("H61").Value = 3845.3447
$ws.Range("I61").Value = 3661.2
$ws.Range("K61").Value = 3661.2
$ws.Range("M61").Value = -3449.2
$ws.Range("H74").Value = 3430.25
$ws.Range("I74").Value = 2117
$ws.Range("J74").Value = 4368.2856
$ws.Range("K74").Value = 2117
$ws.Range("L74").Value = 4368.2856
$ws.Range("M74").Value = -1243
$ws.Range("N74").Value = -6116.2856
$ws.Range("H77").Value = 3430.25
$ws.Range("I77").Value = 2117
$ws.Range("J77").Value = 4368.2856
$ws.Range("K77").Value = 10585
$ws.Range("L77").Value = 21841.428
$ws.Range("M77").Value = -6217
$ws.Range("N77").Value = -30577.428
$ws.Range("H122").Value = 1051573.1
$ws.Range("I122").Value = 3614.8708
$ws.Range("K122").Value = 10844.6124
$ws.Range("M122").Value = -8394.6124
$ws.Range("H135").Value = 349889
$ws.Range("J135").Value = 349889
$ws.Range("L135").Value = 349889
$ws.Range("N135").Value = -360029
$ws.Range("H136").Value = 3845.3447
$ws.Range("I136").Value = 3661.2
$ws.Range("K136").Value = 10983.6
$ws.Range("M136").Value = -8433.599999999999
$ws.Range("H139").Value = 124005
$ws.Range("J139").Value = 124005
$ws.Range("L139").Value = 124005
$ws.Range("N139").Value = -134285

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 20091.572
$ws.Range("I99").Value = 27117.857
$ws.Range("K99").Value = 27117.857
$ws.Range("M99").Value = -25619.857
$ws.Range("H105").Value = 173666.67
$ws.Range("I105").Value = 1000000
$ws.Range("K105").Value = 1000000
$ws.Range("M105").Value = -998253
$ws.Range("H134").Value = 7792.8945
$ws.Range("J134").Value = 888
$ws.Range("L134").Value = 2664
$ws.Range("N134").Value = -7734

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 450
$ws.Range("I10").Value = 450
$ws.Range("K10").Value = 450
$ws.Range("M10").Value = -311
$ws.Range("H58").Value = 2368.125
$ws.Range("I58").Value = 2383.037
$ws.Range("J58").Value = 2287.6
$ws.Range("K58").Value = 2383.037
$ws.Range("L58").Value = 2287.6
$ws.Range("M58").Value = -2180.037
$ws.Range("N58").Value = -2693.6
$ws.Range("H86").Value = 9249.637000000001
$ws.Range("I86").Value = 7917
$ws.Range("J86").Value = 10848.8
$ws.Range("K86").Value = 7917
$ws.Range("L86").Value = 10848.8
$ws.Range("M86").Value = -6794
$ws.Range("N86").Value = -13094.8
$ws.Range("H89").Value = 9249.637000000001
$ws.Range("I89").Value = 7917
$ws.Range("J89").Value = 10848.8
$ws.Range("K89").Value = 39585
$ws.Range("L89").Value = 54244
$ws.Range("M89").Value = -33969
$ws.Range("N89").Value = -65476
$ws.Range("H132").Value = 2992.8
$ws.Range("I132").Value = 2991.923
$ws.Range("J132").Value = 2998.5
$ws.Range("K132").Value = 8975.769
$ws.Range("L132").Value = 8995.5
$ws.Range("M132").Value = -6445.769
$ws.Range("N132").Value = -14055.5
$ws.Range("H134").Value = 4055.1
$ws.Range("I134").Value = 7413.25
$ws.Range("J134").Value = 1816.3334
$ws.Range("K134").Value = 22239.75
$ws.Range("L134").Value = 5449.0002
$ws.Range("M134").Value = -19704.75
$ws.Range("N134").Value = -10519.0002
$ws.Range("H136").Value = 2368.125
$ws.Range("I136").Value = 2383.037
$ws.Range("J136").Value = 2287.6
$ws.Range("K136").Value = 7149.110999999999
$ws.Range("L136").Value = 6862.799999999999
$ws.Range("M136").Value = -4599.110999999999
$ws.Range("N136").Value = -11962.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 425.5
$ws.Range("I23").Value = 448.375
$ws.Range("J23").Value = 414.0625
$ws.Range("K23").Value = 1345.125
$ws.Range("L23").Value = 1242.1875
$ws.Range("M23").Value = -1110.125
$ws.Range("N23").Value = -1712.1875
$ws.Range("H88").Value = 5598
$ws.Range("J88").Value = 5598
$ws.Range("L88").Value = 16794
$ws.Range("N88").Value = -17650
$ws.Range("H91").Value = 5598
$ws.Range("J91").Value = 5598
$ws.Range("L91").Value = 16794
$ws.Range("N91").Value = -19758
$ws.Range("H114").Value = 5499.6665
$ws.Range("I114").Value = 2499
$ws.Range("J114").Value = 7000
$ws.Range("K114").Value = 7497
$ws.Range("L114").Value = 21000
$ws.Range("M114").Value = -4243
$ws.Range("N114").Value = -27508
$ws.Range("H132").Value = 50815.4
$ws.Range("J132").Value = 101013.6
$ws.Range("L132").Value = 909122.4
$ws.Range("N132").Value = -914182.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 12659.182
$ws.Range("I97").Value = 13725.2
$ws.Range("J97").Value = 1999
$ws.Range("K97").Value = 13725.2
$ws.Range("L97").Value = 1999
$ws.Range("M97").Value = -13229.2
$ws.Range("N97").Value = -2991
$ws.Range("H102").Value = 9032
$ws.Range("I102").Value = 10409.733
$ws.Range("K102").Value = 10409.733
$ws.Range("M102").Value = -8787.733
$ws.Range("H113").Value = 8627.764999999999
$ws.Range("I113").Value = 10581
$ws.Range("K113").Value = 10581
$ws.Range("M113").Value = -8411
$ws.Range("H116").Value = 89000
$ws.Range("J116").Value = 89000
$ws.Range("L116").Value = 89000
$ws.Range("N116").Value = -98178

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8771.944
$ws.Range("I16").Value = 9849.714
$ws.Range("K16").Value = 9849.714
$ws.Range("M16").Value = -9679.714
$ws.Range("H119").Value = 19500
$ws.Range("J119").Value = 19500
$ws.Range("L119").Value = 19500
$ws.Range("N119").Value = -29176
$ws.Range("H122").Value = 4779.1763
$ws.Range("I122").Value = 5408.3076
$ws.Range("J122").Value = 2734.5
$ws.Range("K122").Value = 16224.9228
$ws.Range("L122").Value = 8203.5
$ws.Range("M122").Value = -13774.9228
$ws.Range("N122").Value = -13103.5
$ws.Range("H136").Value = 3692.4285
$ws.Range("I136").Value = 2637.4546
$ws.Range("K136").Value = 7912.3638
$ws.Range("M136").Value = -5362.3638

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4104.8887
$ws.Range("I122").Value = 1514.3704
$ws.Range("K122").Value = 4543.1112
$ws.Range("M122").Value = -2093.1112
$ws.Range("H126").Value = 34081.46
$ws.Range("I126").Value = 42066.9
$ws.Range("K126").Value = 126200.7
$ws.Range("M126").Value = -123730.7
$ws.Range("H136").Value = 911709.4399999999
$ws.Range("I136").Value = 1717010.9
$ws.Range("J136").Value = 5745.25
$ws.Range("K136").Value = 5151032.699999999
$ws.Range("L136").Value = 17235.75
$ws.Range("M136").Value = -5148482.699999999
$ws.Range("N136").Value = -22335.75
